$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026490916865734
$ws.Range("D2").Value = 1.026513451475777
$ws.Range("E2").Value = 1.035281354641555
$ws.Range("F2").Value = 1.043026249814855
$ws.Range("I2").Value = 1.028140619206577
$ws.Range("J2").Value = 1.031654064059758
$ws.Range("K2").Value = 1.029335878222202
$ws.Range("L2").Value = 1.03807840888099
$ws.Range("M2").Value = 1.045801277101354
$ws.Range("N2").Value = 1.033119131163119

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027842156397847
$ws.Range("D3").Value = 1.026832767717716
$ws.Range("E3").Value = 1.03652365096261
$ws.Range("F3").Value = 1.044448464980237
$ws.Range("I3").Value = 1.028119731977348
$ws.Range("J3").Value = 1.032643018763664
$ws.Range("K3").Value = 1.029464047931005
$ws.Range("L3").Value = 1.039128916684577
$ws.Range("M3").Value = 1.047032844707646
$ws.Range("N3").Value = 1.034109490296118

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028716173043294
$ws.Range("D4").Value = 1.027038990663314
$ws.Range("E4").Value = 1.037327489195902
$ws.Range("F4").Value = 1.045369012804245
$ws.Range("I4").Value = 1.02810454699596
$ws.Range("J4").Value = 1.033282201577143
$ws.Range("K4").Value = 1.029545830079758
$ws.Range("L4").Value = 1.039808101438636
$ws.Range("M4").Value = 1.047829512859632
$ws.Range("N4").Value = 1.034749580822508

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02908353615165
$ws.Range("D5").Value = 1.027125589626558
$ws.Range("E5").Value = 1.037665423616462
$ws.Range("F5").Value = 1.045756082122717
$ws.Range("I5").Value = 1.028097762610936
$ws.Range("J5").Value = 1.033550740576031
$ws.Range("K5").Value = 1.029579933844576
$ws.Range("L5").Value = 1.040093498731202
$ws.Range("M5").Value = 1.048164378682872
$ws.Range("N5").Value = 1.035018501177576

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029145213854658
$ws.Range("D6").Value = 1.027140124207022
$ws.Range("E6").Value = 1.037722164438539
$ws.Range("F6").Value = 1.045821077047701
$ws.Range("I6").Value = 1.028096599963333
$ws.Range("J6").Value = 1.033595819370691
$ws.Range("K6").Value = 1.029585643700015
$ws.Range("L6").Value = 1.040141410541831
$ws.Range("M6").Value = 1.048220601097659
$ws.Range("N6").Value = 1.035063643989295

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028721082050322
$ws.Range("D7").Value = 1.027040148187686
$ws.Range("E7").Value = 1.037332004684755
$ws.Range("F7").Value = 1.045374184559297
$ws.Range("I7").Value = 1.028104457918005
$ws.Range("J7").Value = 1.033285790486255
$ws.Range("K7").Value = 1.029546286868166
$ws.Range("L7").Value = 1.039811915446027
$ws.Range("M7").Value = 1.047833987558758
$ws.Range("N7").Value = 1.034753174828282

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026947644726116
$ws.Range("D8").Value = 1.026621446480098
$ws.Range("E8").Value = 1.03570119870601
$ws.Range("F8").Value = 1.043506838681758
$ws.Range("I8").Value = 1.02813390540539
$ws.Range("J8").Value = 1.031988440333347
$ws.Range("K8").Value = 1.029379431388346
$ws.Range("L8").Value = 1.038433551539637
$ws.Range("M8").Value = 1.046217542688354
$ws.Range("N8").Value = 1.033453982289349

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023819926531252
$ws.Range("D9").Value = 1.025880711243438
$ws.Range("E9").Value = 1.032827278913259
$ws.Range("F9").Value = 1.04021829376622
$ws.Range("I9").Value = 1.028173049635261
$ws.Range("J9").Value = 1.029696565096052
$ws.Range("K9").Value = 1.029076656696341
$ws.Range("L9").Value = 1.036000241370182
$ws.Range("M9").Value = 1.043367172052731
$ws.Range("N9").Value = 1.031158852326335

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021732687406677
$ws.Range("D10").Value = 1.025385058414673
$ws.Range("E10").Value = 1.030910964937144
$ws.Range("F10").Value = 1.038026988218159
$ws.Range("I10").Value = 1.028190631074268
$ws.Range("J10").Value = 1.028164574574658
$ws.Range("K10").Value = 1.02886901662911
$ws.Range("L10").Value = 1.034374845117876
$ws.Range("M10").Value = 1.041465360200861
$ws.Range("N10").Value = 1.029624686202679

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020828330211808
$ws.Range("D11").Value = 1.025170031100778
$ws.Range("E11").Value = 1.030081045017767
$ws.Range("F11").Value = 1.037078321730686
$ws.Range("I11").Value = 1.028196236062541
$ws.Range("J11").Value = 1.027500197155752
$ws.Range("K11").Value = 1.028777753015332
$ws.Range("L11").Value = 1.033670230717229
$ws.Range("M11").Value = 1.040641432365732
$ws.Range("N11").Value = 1.028959365291633

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020492320724775
$ws.Range("D12").Value = 1.02509010175386
$ws.Range("E12").Value = 1.029772750100411
$ws.Range("F12").Value = 1.036725967244162
$ws.Range("I12").Value = 1.028198017265132
$ws.Range("J12").Value = 1.027253261528089
$ws.Range("K12").Value = 1.028743651841878
$ws.Range("L12").Value = 1.033408380493541
$ws.Range("M12").Value = 1.040335319692193
$ws.Range("N12").Value = 1.028712078987055

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020564400104801
$ws.Range("D13").Value = 1.025107249484421
$ws.Range("E13").Value = 1.029838881637243
$ws.Range("F13").Value = 1.036801547541351
$ws.Range("I13").Value = 1.028197648785465
$ws.Range("J13").Value = 1.027306237190017
$ws.Range("K13").Value = 1.028750975765278
$ws.Range("L13").Value = 1.033464553975997
$ws.Range("M13").Value = 1.040400985102832
$ws.Range("N13").Value = 1.028765129880499

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020800557450148
$ws.Range("D14").Value = 1.025163425298077
$ws.Range("E14").Value = 1.030055561812355
$ws.Range("F14").Value = 1.037049195569475
$ws.Range("I14").Value = 1.028196389426981
$ws.Range("J14").Value = 1.027479788589222
$ws.Range("K14").Value = 1.028774938309027
$ws.Range("L14").Value = 1.033648588669355
$ws.Range("M14").Value = 1.04061613043405
$ws.Range("N14").Value = 1.028938927742598

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020946049567315
$ws.Range("D15").Value = 1.025198029374264
$ws.Range("E15").Value = 1.030189062087951
$ws.Range("F15").Value = 1.037201782521287
$ws.Range("I15").Value = 1.028195573671256
$ws.Range("J15").Value = 1.027586698515009
$ws.Range("K15").Value = 1.028789675730796
$ws.Range("L15").Value = 1.033761961831765
$ws.Range("M15").Value = 1.040748679270769
$ws.Range("N15").Value = 1.029045989492742

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021792694161462
$ws.Range("D16").Value = 1.025399320714895
$ws.Range("E16").Value = 1.030966040613384
$ws.Range("F16").Value = 1.038089951519511
$ws.Range("I16").Value = 1.028190216864643
$ws.Range("J16").Value = 1.028208645405253
$ws.Range("K16").Value = 1.028875045105079
$ws.Range("L16").Value = 1.034421590726031
$ws.Range("M16").Value = 1.041520032120772
$ws.Range("N16").Value = 1.029668819618909

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022323615771931
$ws.Range("D17").Value = 1.025525478374841
$ws.Range("E17").Value = 1.031453377401597
$ws.Range("F17").Value = 1.038647121976612
$ws.Range("I17").Value = 1.028186319765488
$ws.Range("J17").Value = 1.028598501782261
$ws.Range("K17").Value = 1.02892823368536
$ws.Range("L17").Value = 1.034835139739719
$ws.Range("M17").Value = 1.042003762906309
$ws.Range("N17").Value = 1.030059229636689

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022633238848842
$ws.Range("D18").Value = 1.025599024633332
$ws.Range("E18").Value = 1.03173761916036
$ws.Range("F18").Value = 1.038972128104004
$ws.Range("I18").Value = 1.02818385285959
$ws.Range("J18").Value = 1.02882580069113
$ws.Range("K18").Value = 1.02895912692504
$ws.Range("L18").Value = 1.035076278268112
$ws.Range("M18").Value = 1.042285873631966
$ws.Range("N18").Value = 1.03028685133608

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02273880323527
$ws.Range("D19").Value = 1.025624095254975
$ws.Range("E19").Value = 1.031834536082859
$ws.Range("F19").Value = 1.039082950004158
$ws.Range("I19").Value = 1.028182978809425
$ws.Range("J19").Value = 1.02890328728206
$ws.Range("K19").Value = 1.028969638492465
$ws.Range("L19").Value = 1.035158487224475
$ws.Range("M19").Value = 1.042382059235768
$ws.Range("N19").Value = 1.030364447966857

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022266658587825
$ws.Range("D20").Value = 1.025511946906694
$ws.Range("E20").Value = 1.031401092208852
$ws.Range("F20").Value = 1.038587341033059
$ws.Range("I20").Value = 1.028186757925733
$ws.Range("J20").Value = 1.02855668401827
$ws.Range("K20").Value = 1.028922540567382
$ws.Range("L20").Value = 1.034790777901108
$ws.Range("M20").Value = 1.041951867472192
$ws.Range("N20").Value = 1.030017352486677

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020731017553108
$ws.Range("D21").Value = 1.025146884520152
$ws.Range("E21").Value = 1.029991755644657
$ws.Range("F21").Value = 1.036976268817576
$ws.Range("I21").Value = 1.028196768571164
$ws.Range("J21").Value = 1.02742868637165
$ws.Range("K21").Value = 1.028767887494854
$ws.Range("L21").Value = 1.033594398522078
$ws.Range("M21").Value = 1.040552777442377
$ws.Range("N21").Value = 1.028887752954016

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019764970093797
$ws.Range("D22").Value = 1.024917017289605
$ws.Range("E22").Value = 1.029105498996336
$ws.Range("F22").Value = 1.035963451356113
$ws.Range("I22").Value = 1.028201323127522
$ws.Range("J22").Value = 1.026718563481539
$ws.Range("K22").Value = 1.028669484195791
$ws.Range("L22").Value = 1.032841461854167
$ws.Range("M22").Value = 1.039672711468772
$ws.Range("N22").Value = 1.028176621607948

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020277142110416
$ws.Range("D23").Value = 1.025038905493158
$ws.Range("E23").Value = 1.029575336107435
$ws.Range("F23").Value = 1.036500354652359
$ws.Range("I23").Value = 1.028199073251875
$ws.Range("J23").Value = 1.027095100184892
$ws.Range("K23").Value = 1.028721759700692
$ws.Range("L23").Value = 1.033240677867888
$ws.Range("M23").Value = 1.040139290747095
$ws.Range("N23").Value = 1.028553693036609

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022292395255737
$ws.Range("D24").Value = 1.025518061316701
$ws.Range("E24").Value = 1.031424717678507
$ws.Range("F24").Value = 1.038614353412237
$ws.Range("I24").Value = 1.028186560538906
$ws.Range("J24").Value = 1.028575579967526
$ws.Range("K24").Value = 1.028925113446604
$ws.Range("L24").Value = 1.034810823346057
$ws.Range("M24").Value = 1.04197531691229
$ws.Range("N24").Value = 1.030036275270348

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024628864977666
$ws.Range("D25").Value = 1.026072541693807
$ws.Range("E25").Value = 1.033570305248743
$ws.Range("F25").Value = 1.041068256048337
$ws.Range("I25").Value = 1.028164434695646
$ws.Range("J25").Value = 1.030289773934296
$ws.Range("K25").Value = 1.029155958733625
$ws.Range("L25").Value = 1.036629858356135
$ws.Range("M25").Value = 1.044104320226557
$ws.Range("N25").Value = 1.031752903589171
